$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet (Through 2022-10-24 -> Through 2022-10-25)
$ws.Name = "Through 2022-10-25"

# Update the "October" label text (cell A11) to reflect the new through-date
$ws.Range("A11").Value = "October (through 10-25)"

# Update October row (row 11) values
$ws.Range("B11").Value = 23
$ws.Range("C11").Value = 41
$ws.Range("D11").Value = 58
$ws.Range("E11").Value = 55
$ws.Range("F11").Value = 43
$ws.Range("G11").Value = 125
$ws.Range("H11").Value = 157
$ws.Range("I11").Value = 94

# Update Total row (row 12) values
$ws.Range("B12").Value = 249
$ws.Range("C12").Value = 470
$ws.Range("D12").Value = 685
$ws.Range("E12").Value = 603
$ws.Range("F12").Value = 465
$ws.Range("G12").Value = 1026
$ws.Range("H12").Value = 1404
$ws.Range("I12").Value = 1371
